$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comp_type_dmg_algo")
$ws.Columns("G").Insert()
